$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 4450.6562
$ws.Range("J17").Value = 4450.6562
$ws.Range("L17").Value = 13351.9686
$ws.Range("N17").Value = -13687.9686
# Row 63
$ws.Range("H63").Value = 33000
$ws.Range("J63").Value = 33000
$ws.Range("L63").Value = 33000
$ws.Range("N63").Value = -34248
# Row 66
$ws.Range("H66").Value = 33000
$ws.Range("J66").Value = 33000
$ws.Range("L66").Value = 99000
$ws.Range("N66").Value = -105240
# Row 113
$ws.Range("H113").Value = 3466.9333
$ws.Range("J113").Value = 2560
$ws.Range("L113").Value = 2560
$ws.Range("N113").Value = -9068
# Row 137
$ws.Range("H137").Value = 3374.5483
$ws.Range("I137").Value = 973.0625
$ws.Range("J137").Value = 4209.8477
$ws.Range("K137").Value = 2919.1875
$ws.Range("L137").Value = 12629.5431
$ws.Range("M137").Value = -369.1875
$ws.Range("N137").Value = -17729.5431
# Row 138
$ws.Range("H138").Value = 1872.262
$ws.Range("I138").Value = 1400.1034
$ws.Range("K138").Value = 4200.3102
$ws.Range("M138").Value = 939.6898000000001
# Row 140
$ws.Range("H140").Value = 39271.285
$ws.Range("J140").Value = 39271.285
$ws.Range("L140").Value = 39271.285
$ws.Range("N140").Value = -49631.285
# Row 141
$ws.Range("H141").Value = 4924.7617
$ws.Range("I141").Value = 2918.2144
$ws.Range("J141").Value = 8937.857
$ws.Range("K141").Value = 8754.643199999999
$ws.Range("L141").Value = 26813.571
$ws.Range("M141").Value = -3574.643199999999
$ws.Range("N141").Value = -37173.571

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 800
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 1000
$ws.Range("M5").Value = -888
# Row 61
$ws.Range("H61").Value = 2338.4138
$ws.Range("I61").Value = 1335.1
$ws.Range("K61").Value = 1335.1
$ws.Range("M61").Value = -1123.1
# Row 74
$ws.Range("H74").Value = 1063.5385
$ws.Range("I74").Value = 740.2222
$ws.Range("J74").Value = 1791
$ws.Range("K74").Value = 740.2222
$ws.Range("L74").Value = 1791
$ws.Range("M74").Value = 133.7778
$ws.Range("N74").Value = -3539
# Row 77
$ws.Range("H77").Value = 1063.5385
$ws.Range("I77").Value = 740.2222
$ws.Range("J77").Value = 1791
$ws.Range("K77").Value = 3701.111
$ws.Range("L77").Value = 8955
$ws.Range("M77").Value = 666.8889999999997
$ws.Range("N77").Value = -17691
# Row 97
$ws.Range("H97").Value = 603.3043
$ws.Range("I97").Value = 361
$ws.Range("J97").Value = 980.2222
$ws.Range("K97").Value = 361
$ws.Range("L97").Value = 980.2222
$ws.Range("M97").Value = 135
$ws.Range("N97").Value = -1972.2222
# Row 122
$ws.Range("H122").Value = 2853.8948
$ws.Range("I122").Value = 3148
$ws.Range("J122").Value = 2216.6667
$ws.Range("K122").Value = 9444
$ws.Range("L122").Value = 6650.000100000001
$ws.Range("M122").Value = -6994
$ws.Range("N122").Value = -11550.0001
# Row 136
$ws.Range("H136").Value = 2338.4138
$ws.Range("I136").Value = 1335.1
$ws.Range("K136").Value = 4005.3
$ws.Range("M136").Value = -1455.3

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 800
$ws.Range("I4").Value = 1000
$ws.Range("K4").Value = 1000
$ws.Range("M4").Value = -885
# Row 134
$ws.Range("H134").Value = 3575.847
$ws.Range("I134").Value = 1537.079
$ws.Range("J134").Value = 4867.067
$ws.Range("K134").Value = 4611.237
$ws.Range("L134").Value = 14601.201
$ws.Range("M134").Value = -2076.237
$ws.Range("N134").Value = -19671.201

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 339638.2
$ws.Range("I31").Value = 3555.8572
$ws.Range("J31").Value = 406854.66
$ws.Range("K31").Value = 3555.8572
$ws.Range("L31").Value = 406854.66
$ws.Range("M31").Value = -3260.8572
$ws.Range("N31").Value = -407444.66
# Row 34
$ws.Range("H34").Value = 339638.2
$ws.Range("I34").Value = 3555.8572
$ws.Range("J34").Value = 406854.66
$ws.Range("K34").Value = 3555.8572
$ws.Range("L34").Value = 406854.66
$ws.Range("M34").Value = -3353.8572
$ws.Range("N34").Value = -407258.66
# Row 99
$ws.Range("H99").Value = 1819.1818
$ws.Range("I99").Value = 1577.75
$ws.Range("J99").Value = 1957.1428
$ws.Range("K99").Value = 1577.75
$ws.Range("L99").Value = 1957.1428
$ws.Range("M99").Value = -79.75
$ws.Range("N99").Value = -4953.1428
# Row 122
$ws.Range("H122").Value = 120875
$ws.Range("I122").Value = 172221.42
$ws.Range("J122").Value = 1066.6666
$ws.Range("K122").Value = 516664.26
$ws.Range("L122").Value = 3199.9998
$ws.Range("M122").Value = -514214.26
$ws.Range("N122").Value = -8099.9998
# Row 126
$ws.Range("H126").Value = 1819.1818
$ws.Range("I126").Value = 1577.75
$ws.Range("J126").Value = 1957.1428
$ws.Range("K126").Value = 4733.25
$ws.Range("L126").Value = 5871.428400000001
$ws.Range("M126").Value = -2263.25
$ws.Range("N126").Value = -10811.4284
# Row 132
$ws.Range("H132").Value = 29757.7
$ws.Range("I132").Value = 1368.9487
$ws.Range("J132").Value = 130408.73
$ws.Range("K132").Value = 4106.8461
$ws.Range("L132").Value = 391226.19
$ws.Range("M132").Value = -1576.8461
$ws.Range("N132").Value = -396286.19

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 5115.4814
$ws.Range("I131").Value = 20625.8
$ws.Range("J131").Value = 1590.409
$ws.Range("K131").Value = 61877.39999999999
$ws.Range("L131").Value = 4771.227000000001
$ws.Range("M131").Value = -56837.39999999999
$ws.Range("N131").Value = -14851.227

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 123
$ws.Range("H123").Value = 12653.5
$ws.Range("J123").Value = 12653.5
$ws.Range("L123").Value = 12653.5
$ws.Range("N123").Value = -17553.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2340.2144
$ws.Range("I7").Value = 1834.2222
$ws.Range("J7").Value = 3251
$ws.Range("K7").Value = 1834.2222
$ws.Range("L7").Value = 3251
$ws.Range("M7").Value = -1722.2222
$ws.Range("N7").Value = -3475
# Row 122
$ws.Range("H122").Value = 2104.889
$ws.Range("I122").Value = 2111.0588
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 6333.176399999999
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -3883.176399999999
$ws.Range("N122").Value = -10900
# Row 126
$ws.Range("H126").Value = 2340.2144
$ws.Range("I126").Value = 1834.2222
$ws.Range("J126").Value = 3251
$ws.Range("K126").Value = 5502.6666
$ws.Range("L126").Value = 9753
$ws.Range("M126").Value = -3032.6666
$ws.Range("N126").Value = -14693
# Row 136
$ws.Range("H136").Value = 3032.1738
$ws.Range("I136").Value = 2490.9412
$ws.Range("K136").Value = 7472.823600000001
$ws.Range("M136").Value = -4922.823600000001
# Row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 30000
$ws.Range("J64").Value = 30000
$ws.Range("L64").Value = 30000
$ws.Range("N64").Value = -30496
# Row 67
$ws.Range("H67").Value = 30000
$ws.Range("J67").Value = 30000
$ws.Range("L67").Value = 30000
$ws.Range("N67").Value = -31716
